$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-23 18:18:20"
$ws.Range("O2").Value = "6.2 °C"
$ws.Range("E3").Value = "2026-02-23 18:18:22"
$ws.Range("O3").Value = "3.8 °C"
$ws.Range("E4").Value = "2026-02-23 18:18:24"
$ws.Range("K4").Value = "14.8 MJ/m2"
$ws.Range("E5").Value = "2026-02-23 18:18:27"
$ws.Range("O5").Value = "4.6 °C"
$ws.Range("E6").Value = "2026-02-23 18:18:29"
$ws.Range("J6").Value = "1024.5 hPa"
$ws.Range("O6").Value = "14.2 °C"
$ws.Range("E7").Value = "2026-02-23 18:18:32"
$ws.Range("J7").Value = "1024.8 hPa"
$ws.Range("E8").Value = "2026-02-23 18:18:34"
$ws.Range("O8").Value = "14.4 °C"
$ws.Range("E9").Value = "2026-02-23 18:18:37"
$ws.Range("E10").Value = "2026-02-23 18:18:39"
$ws.Range("E11").Value = "2026-02-23 18:18:42"
$ws.Range("E12").Value = "2026-02-23 18:18:44"
$ws.Range("O12").Value = "11.0 °C"
$ws.Range("E13").Value = "2026-02-23 18:18:46"
$ws.Range("J13").Value = "1026.9 hPa"
$ws.Range("O13").Value = "7.1 °C"
$ws.Range("E14").Value = "2026-02-23 18:18:49"
$ws.Range("O14").Value = "13.1 °C"
$ws.Range("E15").Value = "2026-02-23 18:18:51"
$ws.Range("E16").Value = "2026-02-23 18:18:53"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "19%"
$ws.Range("E17").Value = "2026-02-23 18:18:56"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "45%"
$ws.Range("E18").Value = "2026-02-23 18:18:58"
$ws.Range("J18").Value = "1025.0 hPa"
$ws.Range("O18").Value = "11.1 °C"
$ws.Range("E19").Value = "2026-02-23 18:19:01"
$ws.Range("E20").Value = "2026-02-23 18:19:03"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "37%"
$ws.Range("E21").Value = "2026-02-23 18:19:06"
$ws.Range("J21").Value = "1025.9 hPa"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-23 18:19:08"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "21%"
$ws.Range("E23").Value = "2026-02-23 18:19:10"
$ws.Range("E24").Value = "2026-02-23 18:19:13"
$ws.Range("J24").Value = "1026.3 hPa"
$ws.Range("O24").Value = "8.7 °C"
$ws.Range("E25").Value = "2026-02-23 18:19:15"
$ws.Range("E26").Value = "2026-02-23 18:19:18"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "48%"
$ws.Range("J26").Value = "1023.8 hPa"
$ws.Range("O26").Value = "10.4 °C"
$ws.Range("E27").Value = "2026-02-23 18:19:20"
$ws.Range("O27").Value = "5.9 °C"
$ws.Range("E28").Value = "2026-02-23 18:19:22"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "67%"
$ws.Range("J28").Value = "1025.0 hPa"
$ws.Range("E29").Value = "2026-02-23 18:19:25"
$ws.Range("O29").Value = "11.0 °C"
$ws.Range("E30").Value = "2026-02-23 18:19:27"
$ws.Range("J30").Value = "1024.6 hPa"
$ws.Range("E31").Value = "2026-02-23 18:19:30"
$ws.Range("E32").Value = "2026-02-23 18:19:32"
$ws.Range("E33").Value = "2026-02-23 18:19:34"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "45%"
$ws.Range("J33").Value = "1025.4 hPa"
$ws.Range("O33").Value = "8.7 °C"
$ws.Range("E34").Value = "2026-02-23 18:19:37"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "39%"
$ws.Range("O34").Value = "4.4 °C"
$ws.Range("E35").Value = "2026-02-23 18:19:39"
$ws.Range("J35").Value = "1024.9 hPa"
$ws.Range("E36").Value = "2026-02-23 18:19:42"
$ws.Range("O36").Value = "13.1 °C"
$ws.Range("E37").Value = "2026-02-23 18:19:44"
$ws.Range("E38").Value = "2026-02-23 18:19:46"
$ws.Range("O38").Value = "12.5 °C"
$ws.Range("E39").Value = "2026-02-23 18:19:49"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "23%"
$ws.Range("E40").Value = "2026-02-23 18:19:51"
$ws.Range("J40").Value = "1026.2 hPa"
$ws.Range("E41").Value = "2026-02-23 18:19:54"
$ws.Range("E42").Value = "2026-02-23 18:19:56"
$ws.Range("E43").Value = "2026-02-23 18:19:59"
$ws.Range("O43").Value = "10.3 °C"
$ws.Range("E44").Value = "2026-02-23 18:20:01"
$ws.Range("N44").Value = "1.0 °C 17:56 TU"
$ws.Range("E45").Value = "2026-02-23 18:20:03"
$ws.Range("J45").Value = "1026.8 hPa"
$ws.Range("E46").Value = "2026-02-23 18:20:06"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "72%"
$ws.Range("J46").Value = "1026.2 hPa"
$ws.Range("O46").Value = "10.3 °C"
